$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 8 (shifts old rows 8-9 down to 9-10, preserving
# their content/formatting exactly as-is).
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with this week's (2023-07-28) "Primera" price data,
# following the same layout as the other rows.
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(8, 3).Value = "Ñuble"
$ws.Cells.Item(8, 4).Value = 45135
$ws.Cells.Item(8, 5).Value = 16
$ws.Cells.Item(8, 6).Value = 100112039
$ws.Cells.Item(8, 7).Value = "Ciboulette"
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 70
$ws.Cells.Item(8, 11).Value = 2500
$ws.Cells.Item(8, 12).Value = 2500
$ws.Cells.Item(8, 13).Value = 2500
$ws.Cells.Item(8, 14).Value = "`$/docena de atados"
$ws.Cells.Item(8, 15).Value = "Región Metropolitana"
$ws.Cells.Item(8, 16).Value = 833
$ws.Cells.Item(8, 17).Value = 3
$ws.Cells.Item(8, 18).Value = "Hortaliza"
